$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for the Jag2-Notch2 LR-pair sheet.
# Maps each changed cell address to its updated numeric value.
$updates = @{
    "G2" = 20.002957
    "H2" = 60.008871
    "I2" = 0.7920860939997775
    "J2" = 0.7920860939997775
    "M2" = 1.174933333333333
    "N2" = 3.5248
    "O2" = 0.01171850713626266
    "P2" = 0.01171850713626266
    "Q2" = 23.50214094453333
    "R2" = 211.5192685008
    "S2" = 0.009282066545070806
    "T2" = 0.009282066545070805
    "G3" = 20.002957
    "H3" = 60.008871
    "I3" = 0.7920860939997775
    "J3" = 0.7920860939997775
    "O3" = 0.2743256641287217
    "P3" = 0.2743256641287218
    "Q3" = 550.1759181513053
    "R3" = 4951.583263361748
    "S3" = 0.2172895437836141
    "T3" = 0.2172895437836141
    "G4" = 20.002957
    "H4" = 60.008871
    "I4" = 0.7920860939997775
    "J4" = 0.7920860939997775
    "M4" = 39.361408
    "N4" = 118.084224
    "O4" = 0.3925813724534833
    "P4" = 0.3925813724534833
    "Q4" = 787.3445516834561
    "R4" = 7086.100965151104
    "S4" = 0.3109582458837515
    "T4" = 0.3109582458837515
    "G5" = 20.002957
    "H5" = 60.008871
    "I5" = 0.7920860939997775
    "J5" = 0.7920860939997775
    "M5" = 32.221985
    "N5" = 96.665955
    "O5" = 0.3213744562815322
    "P5" = 0.3213744562815322
    "Q5" = 644.5349804096448
    "R5" = 5800.814823686805
    "S5" = 0.254556237787341
    "T5" = 0.2545562377873411
    "I6" = 0.04149178396178559
    "J6" = 0.04149178396178559
    "M6" = 1.174933333333333
    "N6" = 3.5248
    "O6" = 0.01171850713626266
    "P6" = 0.01171850713626266
    "Q6" = 1.231110812444444
    "R6" = 11.079997312
    "S6" = 0.0004862217664524529
    "T6" = 0.0004862217664524528
    "I7" = 0.04149178396178559
    "J7" = 0.04149178396178559
    "O7" = 0.2743256641287217
    "P7" = 0.2743256641287218
    "S7" = 0.01138226119120228
    "T7" = 0.01138226119120228
    "I8" = 0.04149178396178559
    "J8" = 0.04149178396178559
    "M8" = 39.361408
    "N8" = 118.084224
    "O8" = 0.3925813724534833
    "P8" = 0.3925813724534833
    "Q8" = 41.24340812117334
    "R8" = 371.19067309056
    "S8" = 0.01628890149326122
    "T8" = 0.01628890149326122
    "I9" = 0.04149178396178559
    "J9" = 0.04149178396178559
    "M9" = 32.221985
    "N9" = 96.665955
    "O9" = 0.3213744562815322
    "P9" = 0.3213744562815322
    "Q9" = 33.76262550946666
    "R9" = 303.8636295852
    "S9" = 0.01333439951086964
    "T9" = 0.01333439951086964
    "G10" = 2.298356333333333
    "H10" = 6.895068999999999
    "I10" = 0.09101134850660582
    "J10" = 0.09101134850660582
    "M10" = 1.174933333333333
    "N10" = 3.5248
    "O10" = 0.01171850713626266
    "P10" = 0.01171850713626266
    "Q10" = 2.700415467911111
    "R10" = 24.3037392112
    "S10" = 0.001066517136955548
    "T10" = 0.001066517136955548
    "G11" = 2.298356333333333
    "H11" = 6.895068999999999
    "I11" = 0.09101134850660582
    "J11" = 0.09101134850660582
    "O11" = 0.2743256641287217
    "P11" = 0.2743256641287218
    "Q11" = 63.21566885988577
    "R11" = 568.9410197389719
    "S11" = 0.02496674862232519
    "T11" = 0.02496674862232519
    "G12" = 2.298356333333333
    "H12" = 6.895068999999999
    "I12" = 0.09101134850660582
    "J12" = 0.09101134850660582
    "M12" = 39.361408
    "N12" = 118.084224
    "O12" = 0.3925813724534833
    "P12" = 0.3925813724534833
    "Q12" = 90.46654136571733
    "R12" = 814.1988722914559
    "S12" = 0.03572936010556559
    "T12" = 0.03572936010556559
    "G13" = 2.298356333333333
    "H13" = 6.895068999999999
    "I13" = 0.09101134850660582
    "J13" = 0.09101134850660582
    "M13" = 32.221985
    "N13" = 96.665955
    "O13" = 0.3213744562815322
    "P13" = 0.3213744562815322
    "Q13" = 74.05760329732165
    "R13" = 666.5184296758949
    "S13" = 0.02924872264175948
    "T13" = 0.02924872264175948
    "G14" = 1.904387
    "H14" = 5.713160999999999
    "I14" = 0.07541077353183102
    "J14" = 0.07541077353183102
    "M14" = 1.174933333333333
    "N14" = 3.5248
    "O14" = 0.01171850713626266
    "P14" = 0.01171850713626266
    "Q14" = 2.237527765866667
    "R14" = 20.1377498928
    "S14" = 0.0008837016877838489
    "T14" = 0.0008837016877838488
    "G15" = 1.904387
    "H15" = 5.713160999999999
    "I15" = 0.07541077353183102
    "J15" = 0.07541077353183102
    "O15" = 0.2743256641287217
    "P15" = 0.2743256641287218
    "Q15" = 52.37964898091866
    "R15" = 471.416840828268
    "S15" = 0.02068711053158018
    "T15" = 0.02068711053158018
    "G16" = 1.904387
    "H16" = 5.713160999999999
    "I16" = 0.07541077353183102
    "J16" = 0.07541077353183102
    "M16" = 39.361408
    "N16" = 118.084224
    "O16" = 0.3925813724534833
    "P16" = 0.3925813724534833
    "Q16" = 74.959353696896
    "R16" = 674.634183272064
    "S16" = 0.02960486497090504
    "T16" = 0.02960486497090504
    "G17" = 1.904387
    "H17" = 5.713160999999999
    "I17" = 0.07541077353183102
    "J17" = 0.07541077353183102
    "M17" = 32.221985
    "N17" = 96.665955
    "O17" = 0.3213744562815322
    "P17" = 0.3213744562815322
    "Q17" = 61.36312934819499
    "R17" = 552.2681641337549
    "S17" = 0.02423509634156195
    "T17" = 0.02423509634156196
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
